$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New inverse-scored response option label (added to shared strings table)
$newLabel = "3=Not at all true,2=Somewhat true, 1=Very true, 0=Definitely true"

# Rows whose "Value Labels" (column E) switch to the new, inverse-scored label
$rows = @(3, 5, 7, 10, 15, 16, 17, 18, 19, 21, 26)

foreach ($r in $rows) {
    $ws.Range("E$r").Value = $newLabel
}

# Row 26's E cell loses its bottom-border/bold formatting in the target file
$ws.Range("E26").Style = "Normal"

# Match the active selection recorded in the saved workbook
[void]$ws.Range("J10").Select()
